$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = '35.872.77'
$cell.Style = "Normal"
$ws.Range("E2").Value = '  -4.78%  '
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = '1.955.09'
$cell.Style = "Normal"
$ws.Range("E3").Value = '  -5.08%  '
$ws.Range("E4").Value = '  +0.12%  '
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '241.37'
$cell.Style = "Normal"
$ws.Range("E5").Value = '  -4.55%  '
$ws.Range("E6").Value = '  -3.85%  '
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = '61.76'
$cell.Style = "Normal"
$ws.Range("E7").Value = '  -7.49%  '
$ws.Range("E8").Value = '  +0.08%  '
$ws.Range("E9").Value = '  -2.80%  '
$ws.Range("E10").Value = '  -5.85%  '
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = '0.0792'
$cell.Style = "Normal"
$ws.Range("E11").Value = '  +4.54%  '
$ws.Range("E12").Value = '  -1.36%  '
$ws.Range("E13").Value = '  -6.70%  '
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = '21.95'
$cell.Style = "Normal"
$ws.Range("E14").Value = '  +5.51%  '
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = '13.93'
$cell.Style = "Normal"
$ws.Range("E15").Value = '  -9.19%  '
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = '2.242.40'
$cell.Style = "Normal"
$ws.Range("E16").Value = '  -4.87%  '
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = '5.39'
$cell.Style = "Normal"
$ws.Range("E17").Value = '  -4.40%  '
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = '1.965.11'
$cell.Style = "Normal"
$ws.Range("E18").Value = '  -4.44%  '
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = '35.744.60'
$cell.Style = "Normal"
$ws.Range("E19").Value = '  -4.84%  '
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = '70.85'
$cell.Style = "Normal"
$ws.Range("E20").Value = '  -3.44%  '
$ws.Range("E21").Value = '  -3.32%  '
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = '237.52'
$cell.Style = "Normal"
$ws.Range("E22").Value = '  -0.19%  '
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = '5.17'
$cell.Style = "Normal"
$ws.Range("E23").Value = '  -3.76%  '
$ws.Range("E24").Value = '  +0.01%  '
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = '2.52'
$cell.Style = "Normal"
$ws.Range("E25").Value = '  -9.14%  '
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = '2.28'
$cell.Style = "Normal"
$ws.Range("E26").Value = '  -2.78%  '
$ws.Range("E27").Value = '  +1.55%  '
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = '159.05'
$cell.Style = "Normal"
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = '19.74'
$cell.Style = "Normal"
$ws.Range("E29").Value = '  -1.18%  '
$ws.Range("E30").Value = '  +13.71%  '
$ws.Range("E31").Value = '  -2.37%  '
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = '4.83'
$cell.Style = "Normal"
$ws.Range("E32").Value = '  -7.79%  '
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = '1.13'
$cell.Style = "Normal"
$ws.Range("E33").Value = '  -7.61%  '
$ws.Range("E34").Value = '  -0.01%  '
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = '4.38'
$cell.Style = "Normal"
$ws.Range("E35").Value = '  -7.95%  '
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = '6.25'
$cell.Style = "Normal"
$ws.Range("E36").Value = '  +3.21%  '
$ws.Range("E37").Value = '  +0.11%  '
$ws.Range("E38").Value = '  -7.62%  '
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = '1.84'
$cell.Style = "Normal"
$ws.Range("E39").Value = '  +1.33%  '
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = '3.09'
$cell.Style = "Normal"
$ws.Range("E40").Value = '  +14.26%  '
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = '0.0982'
$cell.Style = "Normal"
$ws.Range("E41").Value = '  -5.43%  '
$ws.Range("E42").Value = '  -2.04%  '
$ws.Range("E43").Value = '  -4.10%  '
$ws.Range("E44").Value = '  -4.74%  '
$ws.Range("E45").Value = '  -5.40%  '
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = '91.73'
$cell.Style = "Normal"
$ws.Range("E46").Value = '  -4.30%  '
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = '16.07'
$cell.Style = "Normal"
$ws.Range("E47").Value = '  -6.42%  '
$ws.Range("E48").Value = '  -7.78%  '
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = '1.334.68'
$cell.Style = "Normal"
$ws.Range("E49").Value = '  -6.81%  '
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = '2.77'
$cell.Style = "Normal"
$ws.Range("E50").Value = '  -6.08%  '
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = '2.138.36'
$cell.Style = "Normal"
$ws.Range("E51").Value = '  -4.63%  '
